# AutoCommit_6 ноября 2023 г. 16:03:47_SibNout2023
#
# Adds a new "Лаб6" grade column (J) to the gradebook sheet, gives it the
# same header formatting as the neighbouring "Лаб5" header (I3), fills in a
# few scores that were entered for this grading pass, and moves the active
# cell from G4 to the new J4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column header: J3 = "Лаб6" -----------------------------------
# Copy the formatting of the previous header cell (I3, "Лаб5") onto J3
# first, then set J3's own text so the new header cell keeps the bold /
# centered / bordered header style used by the rest of row 3.
$ws.Range("I3").Copy()
$ws.Range("J3").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("J3").Value = "Лаб6"

# --- New score entries in column G and the new column J ---------------
$ws.Range("G10").Value = 5
$ws.Range("J10").Value = 5

# --- Newly graded "Лаб5" (column I) entries ----------------------------
$ws.Range("I11").Value = 5
$ws.Range("I14").Value = 5
$ws.Range("I21").Value = 5
$ws.Range("I26").Value = 5
$ws.Range("I27").Value = 5
$ws.Range("I28").Value = 5

# --- Newly graded "Лаб3-4" (column H) entry -----------------------------
$ws.Range("H31").Value = 5

# --- Move the active selection to the new last column (J4) -------------
$ws.Range("J4").Select()
